# Glossaire.xlsx bug-fix/documentation update:
#  - Feuil1 ("sheet1.xml"): 5 new glossary terms appended (rows 82-86).
#  - Feuil3 ("sheet2.xml"): a "done" marker added in column B for rows 1-3.

$wb = $excel.ActiveWorkbook

# --- Feuil1: append 5 new single-column glossary rows ---
$ws1 = $wb.Worksheets.Item("Feuil1")

$ws1.Cells.Item(82, 1).Value = "Serializer"
$ws1.Cells.Item(83, 1).Value = "Long"
$ws1.Cells.Item(84, 1).Value = "Array"
$ws1.Cells.Item(85, 1).Value = "Bytes"
$ws1.Cells.Item(86, 1).Value = "Çifteli"

# --- Feuil3: mark the first three rows as "done" in column B ---
$ws3 = $wb.Worksheets.Item("Feuil3")

$ws3.Cells.Item(1, 2).Value = "done"
$ws3.Cells.Item(2, 2).Value = "done"
$ws3.Cells.Item(3, 2).Value = "done"

$ws3.Range("B5").Select() | Out-Null

# Reselect Feuil1 (and its new last cell) last, so it remains the active tab.
$ws1.Range("B88").Select() | Out-Null
